$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's B column ("Total des ventes ...") holds =SUM(C+D) formulas
# whose cached results were never refreshed after D/E started holding text
# in this sheet (see row 5, D5="-4 %"). Switch to manual calculation so our
# edit doesn't force a recalculation (which would turn B6 into #VALUE!).
$excel.Calculation = -4135   # xlCalculationManual

# Row 6: replace the numeric "Ventes de chaï préconfectionné (unités)" and
# "Engagement sur les réseaux sociaux (vues)" entries with text readings.
$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "05:17"
